$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# New identifiers replacing the old handback file ids:
#   9617c70a-46f5-46ac-ab30-2f9d41ca007b  ->  2d407cc9-a69b-446e-b381-78b2585093b0
#   c2ca4700-64fb-4643-8339-4a0566cf22e3  ->  ffff144d9d40-661b-4731-9e5c-105ae3dc13bb
# and refreshed xlf hashes / timestamps produced by a new CI run.
# ---------------------------------------------------------------------------

$file2 = "2d407cc9-a69b-446e-b381-78b2585093b0.md"
$file3 = "ffff144d9d40-661b-4731-9e5c-105ae3dc13bb.md"

$xlfZh = "2d407cc9-a69b-446e-b381-78b2585093b0.818061157a4c319d946f6e613c1fab73a8a0bdb8.zh-cn.xlf"
$xlfDe = "2d407cc9-a69b-446e-b381-78b2585093b0.818061157a4c319d946f6e613c1fab73a8a0bdb8.de-de.xlf"

$overviewDate = "2016-08-15 12:57:30"
$zhHandoffDate = "2016-08-15 12:57:24"
$zhHandbackDate = "2016-08-15 12:57:41"
$deHandoffDate = $overviewDate
$deHandbackDate = "2016-08-15 12:57:48"

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$wsOv = $wb.Worksheets.Item("Overview")

$ovLinks = @{}
foreach ($hl in $wsOv.Hyperlinks) {
    $ovLinks[$hl.Range.Address()] = $hl.Address
}

$wsOv.Range("A2").Value = $file2
$wsOv.Range("G2").Value = $overviewDate

$wsOv.Range("A3").Value = $file3
$wsOv.Range("G3").Value = $overviewDate

$wsOv.Hyperlinks.Delete()
$wsOv.Hyperlinks.Add($wsOv.Range("B2"), $ovLinks["`$B`$2"], [Type]::Missing, [Type]::Missing, "e2e\" + $file2)
$wsOv.Hyperlinks.Add($wsOv.Range("B3"), $ovLinks["`$B`$3"], [Type]::Missing, [Type]::Missing, "e2e\" + $file3)

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$zhLinks = @{}
foreach ($hl in $wsZh.Hyperlinks) {
    $zhLinks[$hl.Range.Address()] = $hl.Address
}

$wsZh.Range("A2").Value = $file2
$wsZh.Range("G2").Value = $xlfZh
$wsZh.Range("H2").Value = $zhHandoffDate
$wsZh.Range("I2").Value = $file2
$wsZh.Range("J2").Value = $xlfZh
$wsZh.Range("K2").Value = $zhHandbackDate

$wsZh.Range("A3").Value = $file3
$wsZh.Range("G3").Value = $xlfZh
$wsZh.Range("H3").Value = $zhHandoffDate
$wsZh.Range("I3").Value = $file3
$wsZh.Range("J3").Value = $xlfZh
$wsZh.Range("K3").Value = $zhHandbackDate

$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), $zhLinks["`$A`$2"], [Type]::Missing, [Type]::Missing, $file2)
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), $zhLinks["`$I`$2"], [Type]::Missing, [Type]::Missing, $file2)
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), $zhLinks["`$A`$3"], [Type]::Missing, [Type]::Missing, $file3)
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), $zhLinks["`$I`$3"], [Type]::Missing, [Type]::Missing, $file3)

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$deLinks = @{}
foreach ($hl in $wsDe.Hyperlinks) {
    $deLinks[$hl.Range.Address()] = $hl.Address
}

$wsDe.Range("A2").Value = $file2
$wsDe.Range("G2").Value = $xlfDe
$wsDe.Range("H2").Value = $deHandoffDate
$wsDe.Range("I2").Value = $file2
$wsDe.Range("J2").Value = $xlfDe
$wsDe.Range("K2").Value = $deHandbackDate

$wsDe.Range("A3").Value = $file3
$wsDe.Range("G3").Value = $xlfDe
$wsDe.Range("H3").Value = $deHandoffDate
$wsDe.Range("I3").Value = $file3
$wsDe.Range("J3").Value = $xlfDe
$wsDe.Range("K3").Value = $deHandbackDate

$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), $deLinks["`$A`$2"], [Type]::Missing, [Type]::Missing, $file2)
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), $deLinks["`$I`$2"], [Type]::Missing, [Type]::Missing, $file2)
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), $deLinks["`$A`$3"], [Type]::Missing, [Type]::Missing, $file3)
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), $deLinks["`$I`$3"], [Type]::Missing, [Type]::Missing, $file3)
